$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 6's formatting from row 5 BEFORE writing any values, so that the
# new date cell (E6) inherits the existing date style instead of Excel
# auto-minting a brand new number format for a previously-unstyled cell.
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(6).RowHeight = 135

# New row 6 values, so new shared strings are interned in the same order
# as the target workbook (B6, then G6, then G3/G5).
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Switched to setup for Staging server provided by the Department for Hosting and security audit,                                  (1) Migration of Test DB, with the exact collation and schema as the real one."
$ws.Range("C6").Value = "Ongoing"
$ws.Range("D6").Value = "LDMS"
$ws.Range("E6").Value2 = (Get-Date -Year 2026 -Month 2 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F6").Value = "Ayush M Srivastava"
$ws.Range("G6").Value = "Reporting Section of LDMS Pending, Validation setup and hardening of TMS pending."

# Row 3: update delivery date and add comment text
$ws.Range("E3").Value2 = (Get-Date -Year 2026 -Month 2 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G3").Value = "Validation setup and hardening of TMS pending."
$ws.Rows.Item(3).RowHeight = 75

# Row 5: update delivery date and add comment text (same text as row 3's comment)
$ws.Range("E5").Value2 = (Get-Date -Year 2026 -Month 2 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G5").Value = "Validation setup and hardening of TMS pending."

$ws.Range("G6").Select()
